$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the 2021 column (M) into the new 2022 column (N)
# for every row that has styled cells, then overwrite with the 2022 values.
$ws.Range("M2").Copy($ws.Range("N2")) | Out-Null
$ws.Range("M3").Copy($ws.Range("N3")) | Out-Null
$ws.Range("M4").Copy($ws.Range("N4")) | Out-Null
$ws.Range("M5").Copy($ws.Range("N5")) | Out-Null
$ws.Range("M6").Copy($ws.Range("N6")) | Out-Null
$ws.Range("M7").Copy($ws.Range("N7")) | Out-Null
$ws.Range("M8").Copy($ws.Range("N8")) | Out-Null
$ws.Range("M9").Copy($ws.Range("N9")) | Out-Null
$ws.Range("M10").Copy($ws.Range("N10")) | Out-Null
$excel.CutCopyMode = 0

# New "2022" column values
$ws.Range("N3").Value = 2022
$ws.Range("N4").Value = 1434
$ws.Range("N5").Value = 12822
$ws.Range("N6").Value = 3099
$ws.Range("N7").Value = 9722
$ws.Range("N8").Value = 14424
$ws.Range("N9").Value = 5279
$ws.Range("N10").Value = 9145

# Move the active selection to the new column's top cell
$ws.Range("N2").Select() | Out-Null
